$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 218.25
$ws.Cells.Item(12, 9).Value = 207.66667
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 207.66667
$ws.Cells.Item(12, 12).Value = 250
$ws.Cells.Item(12, 13).Value = -37.66667000000001
$ws.Cells.Item(12, 14).Value = -590
$ws.Cells.Item(32, 8).Value = 3638
$ws.Cells.Item(32, 10).Value = 3513.7058
$ws.Cells.Item(32, 12).Value = 3513.7058
$ws.Cells.Item(32, 14).Value = -4165.7058
$ws.Cells.Item(38, 8).Value = 5596.794
$ws.Cells.Item(38, 10).Value = 8585
$ws.Cells.Item(38, 12).Value = 25755
$ws.Cells.Item(38, 14).Value = -26499
$ws.Cells.Item(43, 8).Value = 4750.25
$ws.Cells.Item(43, 10).Value = 4667.3335
$ws.Cells.Item(43, 12).Value = 4667.3335
$ws.Cells.Item(43, 14).Value = -4805.3335
$ws.Cells.Item(51, 8).Value = 4920.9
$ws.Cells.Item(51, 9).Value = 4909
$ws.Cells.Item(51, 10).Value = 4922.222
$ws.Cells.Item(51, 11).Value = 4909
$ws.Cells.Item(51, 12).Value = 4922.222
$ws.Cells.Item(51, 13).Value = -4425
$ws.Cells.Item(51, 14).Value = -5890.222
$ws.Cells.Item(53, 8).Value = 2278.375
$ws.Cells.Item(53, 9).Value = 2989.389
$ws.Cells.Item(53, 10).Value = 145.33333
$ws.Cells.Item(53, 11).Value = 2989.389
$ws.Cells.Item(53, 12).Value = 145.33333
$ws.Cells.Item(53, 13).Value = -2352.389
$ws.Cells.Item(53, 14).Value = -1419.33333
$ws.Cells.Item(55, 8).Value = 548.2857
$ws.Cells.Item(55, 9).Value = 242.88889
$ws.Cells.Item(55, 10).Value = 1098
$ws.Cells.Item(55, 11).Value = 242.88889
$ws.Cells.Item(55, 12).Value = 1098
$ws.Cells.Item(55, 13).Value = -28.88889
$ws.Cells.Item(55, 14).Value = -1526
$ws.Cells.Item(58, 8).Value = 5389.9546
$ws.Cells.Item(58, 9).Value = 254.44444
$ws.Cells.Item(58, 11).Value = 763.33332
$ws.Cells.Item(58, 13).Value = -613.33332
$ws.Cells.Item(76, 8).Value = 6768.2
$ws.Cells.Item(76, 10).Value = 8754.5
$ws.Cells.Item(76, 12).Value = 8754.5
$ws.Cells.Item(76, 14).Value = -9384.5
$ws.Cells.Item(79, 8).Value = 6768.2
$ws.Cells.Item(79, 10).Value = 8754.5
$ws.Cells.Item(79, 12).Value = 8754.5
$ws.Cells.Item(79, 14).Value = -10938.5
$ws.Cells.Item(98, 8).Value = 2224.027
$ws.Cells.Item(98, 9).Value = 2159.9092
$ws.Cells.Item(98, 10).Value = 2753
$ws.Cells.Item(98, 11).Value = 2159.9092
$ws.Cells.Item(98, 12).Value = 2753
$ws.Cells.Item(98, 13).Value = -661.9092000000001
$ws.Cells.Item(98, 14).Value = -5749
$ws.Cells.Item(100, 8).Value = 6582916
$ws.Cells.Item(100, 9).Value = 3022.8333
$ws.Cells.Item(100, 11).Value = 3022.8333
$ws.Cells.Item(100, 13).Value = -2481.8333
$ws.Cells.Item(101, 8).Value = 1160.75
$ws.Cells.Item(101, 10).Value = 1499.5
$ws.Cells.Item(101, 12).Value = 4498.5
$ws.Cells.Item(101, 14).Value = -7742.5
$ws.Cells.Item(113, 8).Value = 5648.3438
$ws.Cells.Item(113, 9).Value = 4745.2354
$ws.Cells.Item(113, 10).Value = 5975
$ws.Cells.Item(113, 11).Value = 4745.2354
$ws.Cells.Item(113, 12).Value = 5975
$ws.Cells.Item(113, 13).Value = -1491.2354
$ws.Cells.Item(113, 14).Value = -12483
$ws.Cells.Item(116, 8).Value = 4276
$ws.Cells.Item(116, 9).Value = 4276
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 4276
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -834
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2224.027
$ws.Cells.Item(122, 9).Value = 2159.9092
$ws.Cells.Item(122, 10).Value = 2753
$ws.Cells.Item(122, 11).Value = 6479.7276
$ws.Cells.Item(122, 12).Value = 8259
$ws.Cells.Item(122, 13).Value = -4029.7276
$ws.Cells.Item(122, 14).Value = -13159
$ws.Cells.Item(132, 8).Value = 2698.9473
$ws.Cells.Item(132, 9).Value = 2732.353
$ws.Cells.Item(132, 10).Value = 2415
$ws.Cells.Item(132, 11).Value = 8197.059000000001
$ws.Cells.Item(132, 12).Value = 7245
$ws.Cells.Item(132, 13).Value = -5667.059000000001
$ws.Cells.Item(132, 14).Value = -12305
$ws.Cells.Item(137, 8).Value = 2184.1875
$ws.Cells.Item(137, 9).Value = 2275.7036
$ws.Cells.Item(137, 10).Value = 1690
$ws.Cells.Item(137, 11).Value = 6827.110799999999
$ws.Cells.Item(137, 12).Value = 5070
$ws.Cells.Item(137, 13).Value = -4277.110799999999
$ws.Cells.Item(137, 14).Value = -10170
$ws.Cells.Item(138, 8).Value = 2937.7666
$ws.Cells.Item(138, 9).Value = 2041.8438
$ws.Cells.Item(138, 10).Value = 3961.6785
$ws.Cells.Item(138, 11).Value = 6125.5314
$ws.Cells.Item(138, 12).Value = 11885.0355
$ws.Cells.Item(138, 13).Value = -985.5313999999998
$ws.Cells.Item(138, 14).Value = -22165.0355
$ws.Cells.Item(141, 8).Value = 2053.7
$ws.Cells.Item(141, 9).Value = 1149.4445
$ws.Cells.Item(141, 10).Value = 10192
$ws.Cells.Item(141, 11).Value = 3448.3335
$ws.Cells.Item(141, 12).Value = 30576
$ws.Cells.Item(141, 13).Value = 1731.6665
$ws.Cells.Item(141, 14).Value = -40936

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 194
$ws.Cells.Item(4, 9).Value = 191.25
$ws.Cells.Item(4, 11).Value = 191.25
$ws.Cells.Item(4, 13).Value = -75.25
$ws.Cells.Item(5, 8).Value = 52.666668
$ws.Cells.Item(5, 9).Value = 52.666668
$ws.Cells.Item(5, 11).Value = 52.666668
$ws.Cells.Item(5, 13).Value = 59.333332
$ws.Cells.Item(32, 8).Value = 2445.5781
$ws.Cells.Item(32, 9).Value = 1658.2
$ws.Cells.Item(32, 10).Value = 14256.25
$ws.Cells.Item(32, 11).Value = 1658.2
$ws.Cells.Item(32, 12).Value = 14256.25
$ws.Cells.Item(32, 13).Value = -1371.2
$ws.Cells.Item(32, 14).Value = -14830.25
$ws.Cells.Item(61, 8).Value = 3797.6
$ws.Cells.Item(61, 9).Value = 1996.3334
$ws.Cells.Item(61, 11).Value = 1996.3334
$ws.Cells.Item(61, 13).Value = -1784.3334
$ws.Cells.Item(97, 8).Value = 747.63635
$ws.Cells.Item(97, 9).Value = 653
$ws.Cells.Item(97, 11).Value = 653
$ws.Cells.Item(97, 13).Value = -157
$ws.Cells.Item(102, 8).Value = 19810.963
$ws.Cells.Item(102, 9).Value = 1433.381
$ws.Cells.Item(102, 10).Value = 84132.5
$ws.Cells.Item(102, 11).Value = 1433.381
$ws.Cells.Item(102, 12).Value = 84132.5
$ws.Cells.Item(102, 13).Value = 188.6189999999999
$ws.Cells.Item(102, 14).Value = -87376.5
$ws.Cells.Item(109, 8).Value = 53126.57
$ws.Cells.Item(109, 10).Value = 53126.57
$ws.Cells.Item(109, 12).Value = 53126.57
$ws.Cells.Item(109, 14).Value = -55900.57
$ws.Cells.Item(119, 8).Value = 35000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676
$ws.Cells.Item(122, 8).Value = 4348.8887
$ws.Cells.Item(122, 9).Value = 3953.8
$ws.Cells.Item(122, 10).Value = 5477.7144
$ws.Cells.Item(122, 11).Value = 11861.4
$ws.Cells.Item(122, 12).Value = 16433.1432
$ws.Cells.Item(122, 13).Value = -9411.400000000001
$ws.Cells.Item(122, 14).Value = -21333.1432
$ws.Cells.Item(132, 8).Value = 9117.868
$ws.Cells.Item(132, 9).Value = 4643.383
$ws.Cells.Item(132, 11).Value = 13930.149
$ws.Cells.Item(132, 13).Value = -11400.149
$ws.Cells.Item(135, 8).Value = 60000
$ws.Cells.Item(135, 10).Value = 60000
$ws.Cells.Item(135, 12).Value = 60000
$ws.Cells.Item(135, 14).Value = -70140
$ws.Cells.Item(136, 8).Value = 3797.6
$ws.Cells.Item(136, 9).Value = 1996.3334
$ws.Cells.Item(136, 11).Value = 5989.0002
$ws.Cells.Item(136, 13).Value = -3439.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 52.666668
$ws.Cells.Item(4, 9).Value = 52.666668
$ws.Cells.Item(4, 11).Value = 52.666668
$ws.Cells.Item(4, 13).Value = 62.333332
$ws.Cells.Item(5, 8).Value = 1374.75
$ws.Cells.Item(5, 9).Value = 1374.75
$ws.Cells.Item(5, 11).Value = 1374.75
$ws.Cells.Item(5, 13).Value = -1261.75
$ws.Cells.Item(88, 8).Value = 26600
$ws.Cells.Item(88, 10).Value = 26600
$ws.Cells.Item(88, 12).Value = 26600
$ws.Cells.Item(88, 14).Value = -27412
$ws.Cells.Item(91, 8).Value = 26600
$ws.Cells.Item(91, 10).Value = 26600
$ws.Cells.Item(91, 12).Value = 26600
$ws.Cells.Item(91, 14).Value = -29408
$ws.Cells.Item(94, 8).Value = 930.6721
$ws.Cells.Item(94, 9).Value = 695.25
$ws.Cells.Item(94, 11).Value = 695.25
$ws.Cells.Item(94, 13).Value = -244.25
$ws.Cells.Item(99, 8).Value = 43479708
$ws.Cells.Item(99, 9).Value = 62501316
$ws.Cells.Item(99, 11).Value = 62501316
$ws.Cells.Item(99, 13).Value = -62499818
$ws.Cells.Item(105, 8).Value = 2548.1177
$ws.Cells.Item(105, 9).Value = 2548.1177
$ws.Cells.Item(105, 11).Value = 2548.1177
$ws.Cells.Item(105, 13).Value = -801.1176999999998
$ws.Cells.Item(110, 8).Value = 12702
$ws.Cells.Item(110, 10).Value = 12702
$ws.Cells.Item(110, 12).Value = 12702
$ws.Cells.Item(110, 14).Value = -20882
$ws.Cells.Item(134, 8).Value = 4219.654
$ws.Cells.Item(134, 9).Value = 3920.8125
$ws.Cells.Item(134, 10).Value = 4697.8
$ws.Cells.Item(134, 11).Value = 11762.4375
$ws.Cells.Item(134, 12).Value = 14093.4
$ws.Cells.Item(134, 13).Value = -9227.4375
$ws.Cells.Item(134, 14).Value = -19163.4
$ws.Cells.Item(135, 8).Value = 79600
$ws.Cells.Item(135, 10).Value = 79600
$ws.Cells.Item(135, 12).Value = 79600
$ws.Cells.Item(135, 14).Value = -89740

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 6069.2
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 6069.2
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 6069.2
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = -6293.2
$ws.Cells.Item(7, 8).Value = 311.90475
$ws.Cells.Item(7, 10).Value = 46.25
$ws.Cells.Item(7, 12).Value = 46.25
$ws.Cells.Item(7, 14).Value = -272.25
$ws.Cells.Item(16, 8).Value = 6227.364
$ws.Cells.Item(16, 9).Value = 5699.75
$ws.Cells.Item(16, 11).Value = 5699.75
$ws.Cells.Item(16, 13).Value = -5412.75
$ws.Cells.Item(31, 8).Value = 1436.5264
$ws.Cells.Item(31, 9).Value = 1065.25
$ws.Cells.Item(31, 10).Value = 2073
$ws.Cells.Item(31, 11).Value = 1065.25
$ws.Cells.Item(31, 12).Value = 2073
$ws.Cells.Item(31, 13).Value = -770.25
$ws.Cells.Item(31, 14).Value = -2663
$ws.Cells.Item(34, 8).Value = 1436.5264
$ws.Cells.Item(34, 9).Value = 1065.25
$ws.Cells.Item(34, 10).Value = 2073
$ws.Cells.Item(34, 11).Value = 1065.25
$ws.Cells.Item(34, 12).Value = 2073
$ws.Cells.Item(34, 13).Value = -863.25
$ws.Cells.Item(34, 14).Value = -2477
$ws.Cells.Item(43, 8).Value = 61916.668
$ws.Cells.Item(43, 10).Value = 61916.668
$ws.Cells.Item(43, 12).Value = 61916.668
$ws.Cells.Item(43, 14).Value = -62284.668
$ws.Cells.Item(58, 8).Value = 4824.1665
$ws.Cells.Item(58, 9).Value = 4990
$ws.Cells.Item(58, 11).Value = 4990
$ws.Cells.Item(58, 13).Value = -4787
$ws.Cells.Item(93, 8).Value = 9160.4
$ws.Cells.Item(93, 9).Value = 5200.5
$ws.Cells.Item(93, 10).Value = 25000
$ws.Cells.Item(93, 11).Value = 5200.5
$ws.Cells.Item(93, 12).Value = 25000
$ws.Cells.Item(93, 13).Value = -3328.5
$ws.Cells.Item(93, 14).Value = -28744
$ws.Cells.Item(94, 8).Value = 6724.75
$ws.Cells.Item(94, 9).Value = 3449.5
$ws.Cells.Item(94, 11).Value = 3449.5
$ws.Cells.Item(94, 13).Value = -2998.5
$ws.Cells.Item(95, 8).Value = 49999.5
$ws.Cells.Item(95, 10).Value = 49999.5
$ws.Cells.Item(95, 12).Value = 49999.5
$ws.Cells.Item(95, 14).Value = -55491.5
$ws.Cells.Item(96, 8).Value = 4866.6665
$ws.Cells.Item(96, 10).Value = 4866.6665
$ws.Cells.Item(96, 12).Value = 4866.6665
$ws.Cells.Item(96, 14).Value = -10358.6665
$ws.Cells.Item(97, 8).Value = 55184.668
$ws.Cells.Item(97, 10).Value = 54999.5
$ws.Cells.Item(97, 12).Value = 54999.5
$ws.Cells.Item(97, 14).Value = -56981.5
$ws.Cells.Item(101, 8).Value = 61916.668
$ws.Cells.Item(101, 10).Value = 61916.668
$ws.Cells.Item(101, 12).Value = 61916.668
$ws.Cells.Item(101, 14).Value = -68406.66800000001
$ws.Cells.Item(103, 8).Value = 11261
$ws.Cells.Item(103, 9).Value = 11261
$ws.Cells.Item(103, 11).Value = 11261
$ws.Cells.Item(103, 13).Value = -10089
$ws.Cells.Item(107, 8).Value = 870.2909
$ws.Cells.Item(107, 9).Value = 779.4091
$ws.Cells.Item(107, 10).Value = 930.8788
$ws.Cells.Item(107, 11).Value = 779.4091
$ws.Cells.Item(107, 12).Value = 930.8788
$ws.Cells.Item(107, 13).Value = 1140.5909
$ws.Cells.Item(107, 14).Value = -4770.8788
$ws.Cells.Item(113, 8).Value = 6227.364
$ws.Cells.Item(113, 9).Value = 5699.75
$ws.Cells.Item(113, 11).Value = 5699.75
$ws.Cells.Item(113, 13).Value = -3529.75
$ws.Cells.Item(122, 8).Value = 1869.4193
$ws.Cells.Item(122, 9).Value = 1492.1818
$ws.Cells.Item(122, 11).Value = 4476.5454
$ws.Cells.Item(122, 13).Value = -2026.5454
$ws.Cells.Item(132, 8).Value = 1299.9166
$ws.Cells.Item(132, 9).Value = 1235
$ws.Cells.Item(132, 10).Value = 2014
$ws.Cells.Item(132, 11).Value = 3705
$ws.Cells.Item(132, 12).Value = 6042
$ws.Cells.Item(132, 13).Value = -1175
$ws.Cells.Item(132, 14).Value = -11102
$ws.Cells.Item(133, 8).Value = 46749.75
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 46749.75
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 46749.75
$ws.Cells.Item(133, 13).ClearContents()
$ws.Cells.Item(133, 14).Value = -51809.75
$ws.Cells.Item(134, 8).Value = 1997.0834
$ws.Cells.Item(134, 9).Value = 1764.2354
$ws.Cells.Item(134, 10).Value = 2562.5715
$ws.Cells.Item(134, 11).Value = 5292.706200000001
$ws.Cells.Item(134, 12).Value = 7687.7145
$ws.Cells.Item(134, 13).Value = -2757.706200000001
$ws.Cells.Item(134, 14).Value = -12757.7145
$ws.Cells.Item(136, 8).Value = 4824.1665
$ws.Cells.Item(136, 9).Value = 4990
$ws.Cells.Item(136, 11).Value = 14970
$ws.Cells.Item(136, 13).Value = -12420
$ws.Cells.Item(137, 8).Value = 61998
$ws.Cells.Item(137, 10).Value = 61998
$ws.Cells.Item(137, 12).Value = 61998
$ws.Cells.Item(137, 14).Value = -72198
$ws.Cells.Item(138, 8).Value = 100000
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 100000
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 100000
$ws.Cells.Item(138, 13).ClearContents()
$ws.Cells.Item(138, 14).Value = -110280
$ws.Cells.Item(140, 8).Value = 80000
$ws.Cells.Item(140, 10).Value = 80000
$ws.Cells.Item(140, 12).Value = 80000
$ws.Cells.Item(140, 14).Value = -90360
$ws.Cells.Item(141, 8).Value = 424969.25
$ws.Cells.Item(141, 10).Value = 424969.25
$ws.Cells.Item(141, 12).Value = 424969.25
$ws.Cells.Item(141, 14).Value = -435329.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 19607986
$ws.Cells.Item(2, 9).Value = 147.5
$ws.Cells.Item(2, 11).Value = 885
$ws.Cells.Item(2, 13).Value = -772
$ws.Cells.Item(4, 8).Value = 58830736
$ws.Cells.Item(4, 9).Value = 62500160
$ws.Cells.Item(4, 10).Value = 120000
$ws.Cells.Item(4, 11).Value = 187500480
$ws.Cells.Item(4, 12).Value = 360000
$ws.Cells.Item(4, 13).Value = -187500368
$ws.Cells.Item(4, 14).Value = -360224
$ws.Cells.Item(12, 8).Value = 371.55554
$ws.Cells.Item(12, 10).Value = 465.15384
$ws.Cells.Item(12, 12).Value = 1395.46152
$ws.Cells.Item(12, 14).Value = -1741.46152
$ws.Cells.Item(14, 8).Value = 543.8333
$ws.Cells.Item(14, 9).Value = 543.8333
$ws.Cells.Item(14, 11).Value = 1631.4999
$ws.Cells.Item(14, 13).Value = -1458.4999
$ws.Cells.Item(19, 8).Value = 378.33334
$ws.Cells.Item(19, 10).Value = 320
$ws.Cells.Item(19, 12).Value = 960
$ws.Cells.Item(19, 14).Value = -1308
$ws.Cells.Item(23, 8).Value = 9652946
$ws.Cells.Item(23, 10).Value = 10859554
$ws.Cells.Item(23, 12).Value = 32578662
$ws.Cells.Item(23, 14).Value = -32579132
$ws.Cells.Item(39, 8).Value = 6516.909
$ws.Cells.Item(39, 10).Value = 8614.333000000001
$ws.Cells.Item(39, 12).Value = 25842.999
$ws.Cells.Item(39, 14).Value = -26430.999
$ws.Cells.Item(51, 8).Value = 2333
$ws.Cells.Item(51, 9).Value = 2000
$ws.Cells.Item(51, 10).Value = 2499.5
$ws.Cells.Item(51, 11).Value = 6000
$ws.Cells.Item(51, 12).Value = 7498.5
$ws.Cells.Item(51, 13).Value = -5540
$ws.Cells.Item(51, 14).Value = -8418.5
$ws.Cells.Item(54, 8).Value = 9999
$ws.Cells.Item(54, 10).Value = 9999
$ws.Cells.Item(54, 12).Value = 29997
$ws.Cells.Item(54, 14).Value = -31115
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 10289.3
$ws.Cells.Item(81, 9).Value = 7918.3335
$ws.Cells.Item(81, 11).Value = 23755.0005
$ws.Cells.Item(81, 13).Value = -22632.0005
$ws.Cells.Item(84, 8).Value = 10289.3
$ws.Cells.Item(84, 9).Value = 7918.3335
$ws.Cells.Item(84, 11).Value = 71265.0015
$ws.Cells.Item(84, 13).Value = -65649.0015
$ws.Cells.Item(86, 10).Value = 98
$ws.Cells.Item(86, 12).Value = 294
$ws.Cells.Item(86, 14).Value = -2666
$ws.Cells.Item(89, 10).Value = 98
$ws.Cells.Item(89, 12).Value = 882
$ws.Cells.Item(89, 14).Value = -12738
$ws.Cells.Item(92, 8).Value = 1229.6666
$ws.Cells.Item(92, 9).Value = 396.5
$ws.Cells.Item(92, 11).Value = 1189.5
$ws.Cells.Item(92, 13).Value = 58.5
$ws.Cells.Item(97, 8).Value = 268.5
$ws.Cells.Item(97, 9).Value = 212.125
$ws.Cells.Item(97, 11).Value = 636.375
$ws.Cells.Item(97, 13).Value = -140.375
$ws.Cells.Item(99, 8).Value = 5000
$ws.Cells.Item(99, 9).Value = 5000
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 15000
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -12754
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(118, 8).Value = 335999.66
$ws.Cells.Item(118, 10).Value = 335999.66
$ws.Cells.Item(118, 12).Value = 1007998.98
$ws.Cells.Item(118, 14).Value = -1010484.98
$ws.Cells.Item(120, 8).Value = 10499.8
$ws.Cells.Item(122, 8).Value = 1126.5161
$ws.Cells.Item(122, 9).Value = 1836
$ws.Cells.Item(122, 10).Value = 879.73914
$ws.Cells.Item(122, 11).Value = 16524
$ws.Cells.Item(122, 12).Value = 7917.65226
$ws.Cells.Item(122, 13).Value = -14074
$ws.Cells.Item(122, 14).Value = -12817.65226
$ws.Cells.Item(131, 8).Value = 1646.9231
$ws.Cells.Item(131, 10).Value = 1818.3
$ws.Cells.Item(131, 12).Value = 5454.9
$ws.Cells.Item(131, 14).Value = -15534.9
$ws.Cells.Item(132, 8).Value = 1563.25
$ws.Cells.Item(132, 9).Value = 1334.3334
$ws.Cells.Item(132, 10).Value = 2250
$ws.Cells.Item(132, 11).Value = 12009.0006
$ws.Cells.Item(132, 12).Value = 20250
$ws.Cells.Item(132, 13).Value = -9479.000599999999
$ws.Cells.Item(132, 14).Value = -25310
$ws.Cells.Item(137, 8).Value = 3128.375
$ws.Cells.Item(137, 9).Value = 3208.5
$ws.Cells.Item(137, 11).Value = 9625.5
$ws.Cells.Item(137, 13).Value = -4525.5
$ws.Cells.Item(140, 8).Value = 1276.35
$ws.Cells.Item(140, 9).Value = 865.3889
$ws.Cells.Item(140, 10).Value = 4975
$ws.Cells.Item(140, 11).Value = 2596.1667
$ws.Cells.Item(140, 12).Value = 14925
$ws.Cells.Item(140, 13).Value = 2583.8333
$ws.Cells.Item(140, 14).Value = -25285

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 156.33333
$ws.Cells.Item(5, 9).Value = 156.33333
$ws.Cells.Item(5, 11).Value = 156.33333
$ws.Cells.Item(5, 13).Value = -44.33332999999999
$ws.Cells.Item(12, 8).Value = 501
$ws.Cells.Item(12, 9).Value = 168.66667
$ws.Cells.Item(12, 10).Value = 999.5
$ws.Cells.Item(12, 11).Value = 168.66667
$ws.Cells.Item(12, 12).Value = 999.5
$ws.Cells.Item(12, 13).Value = -28.66667000000001
$ws.Cells.Item(12, 14).Value = -1279.5
$ws.Cells.Item(54, 8).Value = 15222
$ws.Cells.Item(54, 10).Value = 15222
$ws.Cells.Item(54, 12).Value = 15222
$ws.Cells.Item(54, 14).Value = -16002
$ws.Cells.Item(92, 8).Value = 36666
$ws.Cells.Item(92, 10).Value = 49999
$ws.Cells.Item(92, 12).Value = 49999
$ws.Cells.Item(92, 14).Value = -53743
$ws.Cells.Item(102, 8).Value = 2923
$ws.Cells.Item(102, 9).Value = 2046.8889
$ws.Cells.Item(102, 10).Value = 4049.4285
$ws.Cells.Item(102, 11).Value = 2046.8889
$ws.Cells.Item(102, 12).Value = 4049.4285
$ws.Cells.Item(102, 13).Value = -424.8888999999999
$ws.Cells.Item(102, 14).Value = -7293.4285
$ws.Cells.Item(107, 8).Value = 1349.7368
$ws.Cells.Item(107, 10).Value = 661.6
$ws.Cells.Item(107, 12).Value = 661.6
$ws.Cells.Item(107, 14).Value = -4501.6
$ws.Cells.Item(113, 8).Value = 981.6667
$ws.Cells.Item(113, 9).Value = 981.6667
$ws.Cells.Item(113, 11).Value = 981.6667
$ws.Cells.Item(113, 13).Value = 1188.3333
$ws.Cells.Item(132, 8).Value = 3338.606
$ws.Cells.Item(132, 9).Value = 2992.037
$ws.Cells.Item(132, 10).Value = 4898.1665
$ws.Cells.Item(132, 11).Value = 8976.110999999999
$ws.Cells.Item(132, 12).Value = 14694.4995
$ws.Cells.Item(132, 13).Value = -6446.110999999999
$ws.Cells.Item(132, 14).Value = -19754.4995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1475
$ws.Cells.Item(16, 10).Value = 1823
$ws.Cells.Item(16, 12).Value = 1823
$ws.Cells.Item(16, 14).Value = -2163
$ws.Cells.Item(46, 8).Value = 1726.5714
$ws.Cells.Item(46, 9).Value = 1295
$ws.Cells.Item(46, 10).Value = 1899.2
$ws.Cells.Item(46, 11).Value = 1295
$ws.Cells.Item(46, 12).Value = 1899.2
$ws.Cells.Item(46, 13).Value = -1107
$ws.Cells.Item(46, 14).Value = -2275.2
$ws.Cells.Item(68, 8).Value = 3912.2666
$ws.Cells.Item(68, 9).Value = 2033.3636
$ws.Cells.Item(68, 11).Value = 2033.3636
$ws.Cells.Item(68, 13).Value = -1284.3636
$ws.Cells.Item(71, 8).Value = 3912.2666
$ws.Cells.Item(71, 9).Value = 2033.3636
$ws.Cells.Item(71, 11).Value = 10166.818
$ws.Cells.Item(71, 13).Value = -6422.817999999999
$ws.Cells.Item(93, 8).Value = 6779.3335
$ws.Cells.Item(93, 9).Value = 6607.6
$ws.Cells.Item(93, 11).Value = 6607.6
$ws.Cells.Item(93, 13).Value = -5359.6
$ws.Cells.Item(110, 8).Value = 20322
$ws.Cells.Item(110, 10).Value = 20322
$ws.Cells.Item(110, 12).Value = 20322
$ws.Cells.Item(110, 14).Value = -28502
$ws.Cells.Item(122, 8).Value = 6141.7144
$ws.Cells.Item(122, 9).Value = 5269.8887
$ws.Cells.Item(122, 10).Value = 7711
$ws.Cells.Item(122, 11).Value = 15809.6661
$ws.Cells.Item(122, 12).Value = 23133
$ws.Cells.Item(122, 13).Value = -13359.6661
$ws.Cells.Item(122, 14).Value = -28033
$ws.Cells.Item(132, 8).Value = 4344.1
$ws.Cells.Item(132, 9).Value = 3680.3125
$ws.Cells.Item(132, 10).Value = 6999.25
$ws.Cells.Item(132, 11).Value = 11040.9375
$ws.Cells.Item(132, 12).Value = 20997.75
$ws.Cells.Item(132, 13).Value = -8510.9375
$ws.Cells.Item(132, 14).Value = -26057.75
$ws.Cells.Item(136, 8).Value = 4952.147
$ws.Cells.Item(136, 9).Value = 4553.793
$ws.Cells.Item(136, 11).Value = 13661.379
$ws.Cells.Item(136, 13).Value = -11111.379

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 25000
$ws.Cells.Item(2, 10).Value = 25000
$ws.Cells.Item(2, 12).Value = 25000
$ws.Cells.Item(2, 14).Value = -25224
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).ClearContents()
$ws.Cells.Item(17, 8).Value = 178187.25
$ws.Cells.Item(17, 9).Value = 236083
$ws.Cells.Item(17, 10).Value = 4500
$ws.Cells.Item(17, 11).Value = 236083
$ws.Cells.Item(17, 12).Value = 4500
$ws.Cells.Item(17, 13).Value = -235911
$ws.Cells.Item(17, 14).Value = -4844
$ws.Cells.Item(18, 8).Value = 23556
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 23556
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 23556
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(18, 14).Value = -23902
$ws.Cells.Item(94, 8).Value = 45000
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 45000
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 45000
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -46802
$ws.Cells.Item(107, 8).Value = 403.21738
$ws.Cells.Item(107, 9).Value = 382.4375
$ws.Cells.Item(107, 11).Value = 1147.3125
$ws.Cells.Item(107, 13).Value = 772.6875
$ws.Cells.Item(119, 8).Value = 35000
$ws.Cells.Item(119, 9).Value = 20000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 11).Value = 20000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 13).Value = -15162
$ws.Cells.Item(119, 14).Value = -59676
$ws.Cells.Item(122, 8).Value = 3999.26
$ws.Cells.Item(122, 9).Value = 4243.6553
$ws.Cells.Item(122, 11).Value = 12730.9659
$ws.Cells.Item(122, 13).Value = -10280.9659
$ws.Cells.Item(132, 8).Value = 1974.409
$ws.Cells.Item(132, 9).Value = 1297.3334
$ws.Cells.Item(132, 11).Value = 3892.0002
$ws.Cells.Item(132, 13).Value = -1362.0002
$ws.Cells.Item(136, 8).Value = 2887.8125
$ws.Cells.Item(136, 9).Value = 2501.0571
$ws.Cells.Item(136, 10).Value = 3929.077
$ws.Cells.Item(136, 11).Value = 7503.1713
$ws.Cells.Item(136, 12).Value = 11787.231
$ws.Cells.Item(136, 13).Value = -4953.1713
$ws.Cells.Item(136, 14).Value = -16887.231
